# Fixed #253 Moving from POI 3.16 to 3.17.
#
# The expected-generation document embeds a Java stack trace (captured verbatim
# as plain text inside a single run). Upgrading POI from 3.16 to 3.17 shifted
# line numbers in M2Doc's own source files and also collapsed one stack frame
# (the separate "caseTemplate" frame disappeared because of an internal
# refactor), so the recorded trace needs to be updated to match the new
# output.

$d = $word.ActiveDocument

# Unique anchors identifying the first and last line of the stack-frame
# block that changed.
$startMarker = "at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:516)"
$endMarker = "at sun.reflect.GeneratedMethodAccessor76.invoke(Unknown Source)"

$rStart = $d.Content
$foundStart = $rStart.Find.Execute($startMarker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundStart) {
    throw "Could not find start marker for stack trace update"
}

$rEnd = $d.Content
$foundEnd = $rEnd.Find.Execute($endMarker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundEnd) {
    throw "Could not find end marker for stack trace update"
}

# New content for the block, one array entry per stack-trace line (each line
# in the original text is separated by a bare LF followed by a TAB, matching
# the verbatim formatting already used for this trace).
$newLines = @(
    "at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:540)",
    "at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseQuery(M2DocEvaluator.java:1)",
    "at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:186)",
    "at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)",
    "at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)",
    "at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)",
    "at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)",
    "at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)",
    "at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)",
    "at org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)",
    "at org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)",
    "at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)",
    "at org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)",
    "at org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)",
    "at org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)",
    "at org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)",
    "at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)",
    "at org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)",
    "at sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)"
)

$lineSeparator = [string]([char]10) + [string]([char]9)
$newText = [string]::Join($lineSeparator, $newLines)

# Replace the whole old block (start of first changed line through end of
# last changed line) with the rebuilt text in one shot, preserving the
# run's existing formatting (bold, red) since we only replace the Range's
# text.
$target = $d.Range($rStart.Start, $rEnd.End)
$target.Text = $newText

Write-Host "Stack trace block updated."
